$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values need to be swapped between row 26 and row 27
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell26 = $ws.Range($col + "26")
    $cell27 = $ws.Range($col + "27")

    $val26 = $cell26.Value()
    $val27 = $cell27.Value()

    $cell26.Value = $val27
    $cell27.Value = $val26
}
